# Include Emissions from Imported Electricity
# Flip the BIEfIE boolean control-lever from "0" (excluded) to "1" (included).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("BIEfIE")

# Core semantic edit: B2 (the lever value, next to its "Boolean" label in A2)
# goes from 0 to 1 so imported-electricity emissions are included in totals.
$ws.Range("B2").Value = 1

# Leave the same selection state behind that the authored workbook has:
# cell B3 selected/active on the BIEfIE sheet ...
$ws.Range("B3").Select()

# ... while keeping the "About" sheet as the active/selected tab, matching
# the saved workbook view state.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
